$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testen F#")
$ws.Activate()

# --- New content in rows 63-69 ---
# Shared strings must be created in the same order the original author typed
# them so the sharedStrings table ends up in the same order as the target:
#   49 -> A66 text, 50 -> A63 text, 51 -> A65 text, 52 -> A64 text
$euro = [char]0x20AC
$ws.Cells.Item(66, 1).Value = "        /// TW(0.05;2;-1000; 0; 0) = $euro 1102.50. [0 = Postnumerando]"
$ws.Cells.Item(63, 1).Value = "Test de berekening van de toekomstige waarde met enkel een hoofdsom"
$ws.Cells.Item(65, 1).Value = "In Excel is dit :   "
$ws.Cells.Item(64, 1).Value = "achteraf, zonder peridieke betalingen."

# Row 67: column headers (re-use existing shared strings: hoofdsom / rente / looptijd)
$ws.Cells.Item(67, 1).Value = "hoofdsom"
$ws.Cells.Item(67, 2).Value = "rente"
$ws.Cells.Item(67, 3).Value = "looptijd"

# Row 68: input values - hoofdsom, rente (as percentage), looptijd
$ws.Cells.Item(68, 1).Value = -1000
$ws.Cells.Item(68, 2).NumberFormat = "0.00%"
$ws.Cells.Item(68, 2).Value = 0.05
$ws.Cells.Item(68, 3).Value = 2

# Row 69: TW (FV) formula, formatted like the other euro-currency cells
$eurFmt = $ws.Cells.Item(18, 2).NumberFormat
$ws.Cells.Item(69, 1).NumberFormat = $eurFmt
$ws.Cells.Item(69, 1).Formula = "=FV(B68, C68, 0, A68, 0)"

# --- Sheet view: scroll + selection moved down towards the new rows ---
$ws.Range("A66:XFD66").Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "edit complete"
